$d = $word.ActiveDocument

# 1. Merge the five adjacent "{ expression.type := ... }" runs (all sharing the
#    same DC2300/22pt formatting) back into a single run with the concatenated
#    text, by replacing the full text with itself.
$oldText1 = "{ expression.type := if expression1 := void then simple_expression.type, else expression1.type }"
$d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $oldText1, 2) | Out-Null

# 2. Remove the stray "{not done: function calls}" run entirely (blue,
#    0000FF) that trails "factor1" in the `factor -> id factor1` production.
$oldText2 = "{not done: function calls}"
$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 3. Flip OverflowPunct off on the "Normal" style's paragraph format.
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.OverflowPunct = $false
